# Translations workbook: add a new tooltip entry for "SetActiveAtStartToggle"
# (ID / ENGLISH / SPANISH columns) as row 48, merge columns B:C into a wider
# wrapped column, and select the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 48 content ------------------------------------------------
# Column order in the sheet is A=ID, B=ENGLISH, C=SPANISH, so the shared
# strings table picks up the new unique strings in ID -> SPANISH -> ENGLISH
# order when we write the cells in that sequence.
$ws.Range("A48").Value = "tooltip.SetActiveAtStartToggle"
$ws.Range("C48").Value = "Define si el objeto estará [00FFFF]activado[-] o [FF0000]no[-] al inicio del nivel."
$ws.Range("B48").Value = "Defines if the object will be [00FFFF]enabled[-] or [FF0000]not[-] at the beginning of the level."

# --- Row heights ---------------------------------------------------------
$ws.Rows.Item(39).RowHeight = 15
$ws.Rows.Item(48).RowHeight = 45

# --- Cell alignment / wrapping for the new row ---------------------------
$a48 = $ws.Range("A48")
$a48.VerticalAlignment = -4160   # xlVAlignTop

$b48 = $ws.Range("B48")
$b48.HorizontalAlignment = -4131 # xlVAlignLeft
$b48.VerticalAlignment = -4160   # xlVAlignTop
$b48.WrapText = $true

$c48 = $ws.Range("C48")
$c48.VerticalAlignment = -4160   # xlVAlignTop
$c48.WrapText = $true

# --- Column widths: merge B & C into one wider, wrapped column width -----
# (45.7109375 chars is the exact target; the host snaps ColumnWidth to the
# nearest 1/6-character increment, so 44.8 is the closest input that lands
# on the nearest achievable stored width.)
$ws.Range("B1:C1").ColumnWidth = 44.8

# --- Selection -------------------------------------------------------------
$ws.Range("B50").Select()

Write-Host "Applied tooltip.SetActiveAtStartToggle row"
